$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: City D
$ws.Range("A5").Value = "City D"
$ws.Range("B5").Value = 1023
$ws.Range("C5").Value = 2025
$ws.Range("D5").Value = 400
$ws.Range("E5").Value = 1034
$ws.Range("F5").Value = 505

# Row 6: City E
$ws.Range("A6").Value = "City E"
$ws.Range("B6").Value = 1600
$ws.Range("C6").Value = 1800
$ws.Range("D6").Value = 1700
$ws.Range("E6").Value = 1700
$ws.Range("F6").Value = 1609

# Row 7: City Ties
$ws.Range("A7").Value = "City Ties"
$ws.Range("B7").Value = 300
$ws.Range("C7").Value = 400
$ws.Range("D7").Value = 300
$ws.Range("E7").Value = 400
$ws.Range("F7").Value = 400

# Row 8: City F
$ws.Range("A8").Value = "City F"
$ws.Range("B8").Value = 456
$ws.Range("C8").Value = 723
$ws.Range("D8").Value = 244
$ws.Range("E8").Value = 500
$ws.Range("F8").Value = 456

$ws.Range("A9").Select() | Out-Null
